# Actualización automática del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with resultado / profit values ---
$ws.Range("G5").Value = "Acierto"
$ws.Range("H5").Value = 0.67

$ws.Range("G28").Value = "Fallo"
$ws.Range("H28").Value = -1

$ws.Range("G45").Value = "Acierto"
$ws.Range("H45").Value = 1.38

$ws.Range("G46").Value = "Fallo"
$ws.Range("H46").Value = -1

# --- Append new rows (54-68) with new match data ---
$newRows = @(
    @(14633480, "2025-09-09", "Jazmin Ortenzi", "Berfu Cengiz", "Gana Berfu Cengiz", 2),
    @(14633477, "2025-09-09", "Valeriya Strakhova", "Ana Candiotto", "Gana Ana Candiotto", 4.33),
    @(14643798, "2025-09-09", "Santiago Rodriguez Taverna", "Rudolf Molleker", "Gana Santiago Rodriguez Taverna", 1.83),
    @(14633450, "2025-09-09", "Andrea Pellegrino", "Mateusz Lange", "Gana Mateusz Lange", 51),
    @(14601373, "2025-09-10", "Facundo Diaz Acosta", "Hynek Barton", "Gana Hynek Barton", 2.38),
    @(14644103, "2025-09-09", "Maxime Janvier", "Calvin Hemery", "Gana Maxime Janvier", 2.38),
    @(14644104, "2025-09-09", "Kenny De Schepper", "Stan Wawrinka", "Gana Kenny De Schepper", 8),
    @(14644855, "2025-09-10", "Cannon Kingsley", "Sean Cuenin", "Gana Cannon Kingsley", 1.53),
    @(14644854, "2025-09-10", "Michael Geerts", "Robin Catry", "Gana Robin Catry", 2.63),
    @(14634123, "2025-09-09", "Aidan Kim", "Ryuki Matsuda", "Gana Aidan Kim", 1.53),
    @(14634114, "2025-09-09", "Alex Rybakov", "Thai-Son Kwiatkowski", "Gana Thai-Son Kwiatkowski", 3.4),
    @(14645805, "2025-09-09", "Edward Winter", "Alfredo Perez", "Gana Edward Winter", 1.83),
    @(14645819, "2025-09-09", "Cooper Williams", "Murphy Cassone", "Gana Cooper Williams", 3.4),
    @(14645817, "2025-09-09", "Rafael Jodar", "Gavin Goode", "Gana Gavin Goode", 11),
    @(14645820, "2025-09-09", "Quinn Vandecasteele", "James Trotter", "Gana Quinn Vandecasteele", 3.75)
)

$startRow = 54
$endRow = $startRow + $newRows.Count - 1

# Column B holds dates stored as plain text (e.g. "2025-09-09"), not Excel
# date serials. Force text formatting before assignment so the COM layer
# doesn't auto-convert the recognizable date strings, then strip the
# number-format back off so the cells end up with no explicit style
# (matching the rest of the data rows).
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

$ws.Range("B$startRow`:B$endRow").ClearFormats()
